$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 54.86376272656823

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
